# MorganPatrone2006a__C_Stationarygenerator_alpha_zero.xlsx
# Commit: "expermits todos no convexos menos el 5to"
#
# Updates the generator's numeric experiment data on several sheets.
# Values that look numeric but must be preserved as literal text (to
# match the original authoring, which stored these as shared strings,
# not spreadsheet numbers) are entered with a leading apostrophe so
# Excel stores them as text instead of coercing them into numbers.
#
# NOTE: the workbook has two sheets whose names differ only by case
# ("Vector_bf" and "Vector_BF"); Worksheets.Item(<name>) resolves
# case-insensitively and would hit the wrong one, so every sheet below
# is addressed by its 1-based tab index instead of by name.

$wb = $excel.ActiveWorkbook

# --- Restricciones_del_lider (sheet 2) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2").Value = "4.5 - x"
$ws2.Range("B2").Value = "'-5.0"
$ws2.Range("D2").Value = "'0.34"
$ws2.Range("A3").Value = "-4.5 + x"
$ws2.Range("B3").Value = "'4.0"
$ws2.Range("D3").Value = "'0.0"

# --- Restricciones_del_follower (sheet 3) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2").Value = "2.8 - y"
$ws3.Range("B2").Value = "'-3.8"
$ws3.Range("D2").Value = "'0.09"
$ws3.Range("E2").Value = "'0.5"
$ws3.Range("F2").Value = "'1.7000000000000002"
$ws3.Range("A3").Value = "-2.8 + y"
$ws3.Range("B3").Value = "'1.7999999999999998"
$ws3.Range("D3").Value = "'0.82"
$ws3.Range("E3").Value = "'-9.200000000000001"
$ws3.Range("F3").Value = "'-0.1"

# --- Punto_modificado (sheet 4) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A2").Value = "'4.5"
$ws4.Range("B2").Value = "'2.8"

# --- Vector_bf (sheet 5) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("A2").Value = "'-5.23"

# --- Vector_BF (sheet 6) ---
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("A2").Value = "'1.34"
$ws6.Range("A3").Value = "'10.700000000000001"
